# Applies the "Preprocessing module scripts from other branch" edit to
# db_defaultEnergyAssets.xlsx
#
# Sheet1 = consumptionAssets
# Sheet2 = productionAssets
# Sheet3 = conversionAssets
# Sheet4 = storageAssets

$wb = $excel.ActiveWorkbook

$wsConsumption = $wb.Worksheets.Item(1)
$wsProduction  = $wb.Worksheets.Item(2)
$wsConversion  = $wb.Worksheets.Item(3)
$wsStorage     = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# productionAssets (sheet2): insert a new "0 MW" solar field row right
# before the existing "Solarpanels_1MW" row, then append several new
# solar field / building solar-panel rows after the existing data.
# ---------------------------------------------------------------------

# Shift the existing row 8 (Solarpanels_1MW) down to row 9 and insert a
# blank row 8 for the new "Solarpanels_0MW" entry.
$wsProduction.Rows.Item(8).Insert()

$wsProduction.Cells.Item(8, 1).Value = 7
$wsProduction.Cells.Item(8, 2).Value = "Solarpanels_0MW"
$wsProduction.Cells.Item(8, 3).Value = "Solar field 0 MW"
$wsProduction.Cells.Item(8, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(8, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(8, 6).Value = 0
$wsProduction.Cells.Item(8, 7).Value = 0

# Renumber the id of the row that got pushed down (was id 7, now id 8).
$wsProduction.Cells.Item(9, 1).Value = 8

# New rows 10-20.
$wsProduction.Cells.Item(10, 1).Value = 9
$wsProduction.Cells.Item(10, 2).Value = "Solarpanels_2MW"
$wsProduction.Cells.Item(10, 3).Value = "Solar field 2 MW"
$wsProduction.Cells.Item(10, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(10, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(10, 6).Value = 2000
$wsProduction.Cells.Item(10, 7).Value = 0

$wsProduction.Cells.Item(11, 1).Value = 10
$wsProduction.Cells.Item(11, 2).Value = "Solarpanels_3MW"
$wsProduction.Cells.Item(11, 3).Value = "Solar field 3 MW"
$wsProduction.Cells.Item(11, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(11, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(11, 6).Value = 3000
$wsProduction.Cells.Item(11, 7).Value = 0

$wsProduction.Cells.Item(12, 1).Value = 11
$wsProduction.Cells.Item(12, 2).Value = "Solarpanels_4MW"
$wsProduction.Cells.Item(12, 3).Value = "Solar field 4 MW"
$wsProduction.Cells.Item(12, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(12, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(12, 6).Value = 4000

$wsProduction.Cells.Item(13, 1).Value = 12
$wsProduction.Cells.Item(13, 2).Value = "Solarpanels_5MW"
$wsProduction.Cells.Item(13, 3).Value = "Solar field 5 MW"
$wsProduction.Cells.Item(13, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(13, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(13, 6).Value = 5000

$wsProduction.Cells.Item(14, 1).Value = 13
$wsProduction.Cells.Item(14, 2).Value = "Solarpanels_6MW"
$wsProduction.Cells.Item(14, 3).Value = "Solar field 6 MW"
$wsProduction.Cells.Item(14, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(14, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(14, 6).Value = 6000

$wsProduction.Cells.Item(15, 1).Value = 14
$wsProduction.Cells.Item(15, 2).Value = "Building_solarpanels_0kWp"
$wsProduction.Cells.Item(15, 3).Value = "Solar panels for building, 0 kWp"
$wsProduction.Cells.Item(15, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(15, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(15, 6).Value = 0
$wsProduction.Cells.Item(15, 7).Value = 0

$wsProduction.Cells.Item(16, 1).Value = 15
$wsProduction.Cells.Item(16, 2).Value = "Building_solarpanels_1kWp"
$wsProduction.Cells.Item(16, 3).Value = "Solar panels for building, 1 kWp"
$wsProduction.Cells.Item(16, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(16, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(16, 6).Value = 1
$wsProduction.Cells.Item(16, 7).Value = 0

$wsProduction.Cells.Item(17, 1).Value = 16
$wsProduction.Cells.Item(17, 2).Value = "Building_solarpanels_2kWp"
$wsProduction.Cells.Item(17, 3).Value = "Solar panels for building, 2 kWp"
$wsProduction.Cells.Item(17, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(17, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(17, 6).Value = 2
$wsProduction.Cells.Item(17, 7).Value = 0

$wsProduction.Cells.Item(18, 1).Value = 17
$wsProduction.Cells.Item(18, 2).Value = "Building_solarpanels_3kWp"
$wsProduction.Cells.Item(18, 3).Value = "Solar panels for building, 3 kWp"
$wsProduction.Cells.Item(18, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(18, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(18, 6).Value = 3
$wsProduction.Cells.Item(18, 7).Value = 0

$wsProduction.Cells.Item(19, 1).Value = 18
$wsProduction.Cells.Item(19, 2).Value = "Building_solarpanels_4kWp"
$wsProduction.Cells.Item(19, 3).Value = "Solar panels for building, 4 kWp"
$wsProduction.Cells.Item(19, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(19, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(19, 6).Value = 4
$wsProduction.Cells.Item(19, 7).Value = 0

$wsProduction.Cells.Item(20, 1).Value = 19
$wsProduction.Cells.Item(20, 2).Value = "Building_solarpanels_5kWp"
$wsProduction.Cells.Item(20, 3).Value = "Solar panels for building, 5 kWp"
$wsProduction.Cells.Item(20, 4).Value = "PRODUCTION"
$wsProduction.Cells.Item(20, 5).Value = "PHOTOVOLTAIC"
$wsProduction.Cells.Item(20, 6).Value = 5
$wsProduction.Cells.Item(20, 7).Value = 0

# ---------------------------------------------------------------------
# conversionAssets (sheet3): append a building gas burner and a
# building heat pump.
# ---------------------------------------------------------------------

$wsConversion.Cells.Item(13, 1).Value = 12
$wsConversion.Cells.Item(13, 2).Value = "Building_gas_burner_60kW"
$wsConversion.Cells.Item(13, 3).Value = "CONVERSION"
$wsConversion.Cells.Item(13, 4).Value = "GAS_BURNER"
$wsConversion.Cells.Item(13, 5).Value = 0
$wsConversion.Cells.Item(13, 6).Value = 60
$wsConversion.Cells.Item(13, 7).Value = 0.95
$wsConversion.Cells.Item(13, 8).Value = 90

$wsConversion.Cells.Item(14, 1).Value = 13
$wsConversion.Cells.Item(14, 2).Value = "Building_heatpump_20kW"
$wsConversion.Cells.Item(14, 3).Value = "CONVERSION"
$wsConversion.Cells.Item(14, 4).Value = "HEAT_PUMP_AIR"
$wsConversion.Cells.Item(14, 5).Value = 20
$wsConversion.Cells.Item(14, 6).Value = 0
$wsConversion.Cells.Item(14, 7).Value = 0.6
$wsConversion.Cells.Item(14, 8).Value = 60
$wsConversion.Cells.Item(14, 9).Value = "AIR"

# ---------------------------------------------------------------------
# storageAssets (sheet4): update the EHGV row to use the new
# ELECTRIC_HEAVY_GOODS_VEHICLE asset type and a larger capacity, bump
# the round-trip efficiency of the two grid batteries to 1, and drop
# the now-unused energy_consumption_kwhpkm column (O).
# ---------------------------------------------------------------------

# Grid_battery_1MWh: stateOfCharge_r 0.8 -> 1
$wsStorage.Cells.Item(14, 7).Value = 1

# EHGV: energyAssetType ELECTRIC_VEHICLE -> ELECTRIC_HEAVY_GOODS_VEHICLE,
# capacityElectricity_kW 50 -> 110
$wsStorage.Cells.Item(15, 4).Value = "ELECTRIC_HEAVY_GOODS_VEHICLE"
$wsStorage.Cells.Item(15, 5).Value = 110

# Grid_battery_10MWh: stateOfCharge_r 0.8 -> 1
$wsStorage.Cells.Item(16, 7).Value = 1

# Remove the energy_consumption_kwhpkm column entirely.
$wsStorage.Columns.Item(15).Delete()

# Narrow column D now that the ELECTRIC_VEHICLE values no longer need
# to be shown as widely.
$wsStorage.Columns.Item(4).ColumnWidth = 20.33

# ---------------------------------------------------------------------
# View / selection bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------

$wsConsumption.Range("F10").Select()
$wsProduction.Range("B23").Select()
$wsStorage.Range("L17").Select()
$wsConversion.Range("C17").Select()
$wsConversion.Activate()
